$wb = $excel.ActiveWorkbook
$wsTranscriptions = $wb.Worksheets.Item("Transcriptions")
$wsAnnotations = $wb.Worksheets.Item("Annotations")

# Add three new annotation rows (110-112) to the Annotations sheet.
$wsAnnotations.Cells.Item(110, 1).Value = "Golden Bowl, The"
$wsAnnotations.Cells.Item(110, 2).Value = "Literary Work"
$wsAnnotations.Cells.Item(110, 3).Value = "lit-gold"
$wsAnnotations.Cells.Item(110, 4).Value = "../resources/annotations.xml#lit-gold"

$wsAnnotations.Cells.Item(111, 1).Value = "Ambassadors, The"
$wsAnnotations.Cells.Item(111, 2).Value = "Literary Work"
$wsAnnotations.Cells.Item(111, 3).Value = "lit-amba"
$wsAnnotations.Cells.Item(111, 4).Value = "../resources/annotations.xml#lit-amba"

$wsAnnotations.Cells.Item(112, 1).Value = "Edward IV"
$wsAnnotations.Cells.Item(112, 2).Value = "Person"
$wsAnnotations.Cells.Item(112, 3).Value = "psn-edw4"
$wsAnnotations.Cells.Item(112, 4).Value = "../resources/annotations.xml#psn-edw4"
$wsAnnotations.Cells.Item(112, 5).Value = "Is this actually Edward IV??? I'm guessing based on Latimer and the c16th but then again… Check against quotes! "

# Restore the selection on the Transcriptions sheet (no longer the active tab).
$wsTranscriptions.Range("B32").Select() | Out-Null

# Make Annotations the active sheet/tab, scrolled & selected as in the target workbook.
$wsAnnotations.Activate() | Out-Null
$wsAnnotations.Range("G74").Select() | Out-Null
